# New PO forecast model
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Weekly Quantity -> append a new weekly row (row 21)
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A21").Value = 45676.99999999999
$wsWeekly.Range("A21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("B21").Value = 1

# ---------------------------------------------------------------------
# Sheet: Monthly Trend -> append a new monthly row (row 8)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A8").Value = 45688.99999999999
$wsMonthly.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMonthly.Range("B8").Value = 1

# ---------------------------------------------------------------------
# Sheet: PO Forecast -> refresh the forecast model values
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Updated forecast quantities for existing dates (rows 2-20, dates unchanged)
$wsForecast.Range("B2").Value = 15
$wsForecast.Range("B3").Value = 20
$wsForecast.Range("B4").Value = 26
$wsForecast.Range("B5").Value = 31
$wsForecast.Range("B6").Value = 36
$wsForecast.Range("B7").Value = 42
$wsForecast.Range("B8").Value = 47
$wsForecast.Range("B9").Value = 53
$wsForecast.Range("B10").Value = 58
$wsForecast.Range("B11").Value = 63
$wsForecast.Range("B12").Value = 69
$wsForecast.Range("B13").Value = 74
$wsForecast.Range("B14").Value = 80
$wsForecast.Range("B15").Value = 90
$wsForecast.Range("B16").Value = 96
$wsForecast.Range("B17").Value = 101
$wsForecast.Range("B18").Value = 107
$wsForecast.Range("B19").Value = 112
$wsForecast.Range("B20").Value = 117

# Rows 21-28: new dates/values (forecast horizon shifted forward)
$wsForecast.Range("A21").Value = 45676.99999999999
$wsForecast.Range("B21").Value = 144
$wsForecast.Range("A22").Value = 45683.99999999999
$wsForecast.Range("B22").Value = 150
$wsForecast.Range("A23").Value = 45690.99999999999
$wsForecast.Range("B23").Value = 155
$wsForecast.Range("A24").Value = 45697.99999999999
$wsForecast.Range("B24").Value = 160
$wsForecast.Range("A25").Value = 45704.99999999999
$wsForecast.Range("B25").Value = 166
$wsForecast.Range("A26").Value = 45711.99999999999
$wsForecast.Range("B26").Value = 171
$wsForecast.Range("A27").Value = 45718.99999999999
$wsForecast.Range("B27").Value = 177
$wsForecast.Range("A28").Value = 45725.99999999999
$wsForecast.Range("B28").Value = 182

# Row 29: brand-new forecast row
$wsForecast.Range("A29").Value = 45732.99999999999
$wsForecast.Range("A29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B29").Value = 187
